$wb = $excel.ActiveWorkbook

# Sheet1: row 2 gets a brand-new SMS log entry (phone numbers, date/time,
# sequence number) and the AZ column (header + value) switches from the
# "SMSTemplate_Msg" sample text to a new "CurrentTime" template.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F2").Value = "0141935690"
$ws1.Range("N2").Value = "2024-03-06"
$ws1.Range("O2").Value = "02:35:55 PM"
$ws1.Range("P2").Value = "2024-03-06 07:24:53 PM"
$ws1.Range("AC2").Value = "2024-03-06"
$ws1.Range("AE2").Value = "1447058961"
$ws1.Range("AK2").Value = "4"
$ws1.Range("AT2").Value = "9685046387"
$ws1.Range("AX2").Value = "5762296884"
$ws1.Range("AZ1").Value = "CurrentTime"
$ws1.Range("AZ2").Value = "CT. Wed, Mar 06, 2024 at 7:30 PM"
$ws1.Range("AZ2").Style = $ws1.Range("AZ1").Style

# Sheet2: same batch of phone-number / mobile-number columns refreshed.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("F2").Value = "0141935690"
$ws2.Range("AE2").Value = "1447058961"
$ws2.Range("AT2").Value = "9685046387"
$ws2.Range("AX2").Value = "5762296884"

# Sheet3: same refresh.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("F2").Value = "0141935690"
$ws3.Range("AE2").Value = "1447058961"
$ws3.Range("AT2").Value = "9685046387"
$ws3.Range("AX2").Value = "5762296884"

# Sheet4: same refresh.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("F2").Value = "0141935690"
$ws4.Range("AE2").Value = "1447058961"
$ws4.Range("AT2").Value = "9685046387"
$ws4.Range("AX2").Value = "5762296884"
